# "Grosse mise à jour" - Répartition.xlsx
#
# 1) Mark several tasks as finished ("✓") or in-progress ("∼") in the
#    "Etat" column (J) of each of the four "étapes" tables.
# 2) Flag two tasks as postponed by writing "Reporté" (white text on a
#    red fill) into column K next to them, and widen column K so the
#    label is readable (column L - a hidden helper column - shrinks
#    to a sliver).
# 3) Move the active selection down to the newly annotated cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "Etat" column (J) updates.
#
# Use an already-correctly-styled cell as the template so the new cells
# pick up the very same number format / fill / font / border that the
# rest of the column already uses, then overwrite the displayed glyph.
# ---------------------------------------------------------------------

$checkTemplate = $ws.Range("J3")   # existing "✓" cell (fontId3/fillId19 style)
$tildeTemplate = $ws.Range("Q5")   # existing "∼" cell (fontId4/fillId20 style) - untouched by this edit

$toCheck = @("J4", "J6", "J25", "J35", "J36", "J37", "J45", "J46")
foreach ($ref in $toCheck) {
    $checkTemplate.Copy()
    $ws.Range($ref).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range($ref).Value2 = "✓"
}

$toTilde = @("J47")
foreach ($ref in $toTilde) {
    $tildeTemplate.Copy()
    $ws.Range($ref).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range($ref).Value2 = "∼"
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) "Reporté" (postponed) markers in column K.
#
# Base the formatting on a neighbouring bordered cell (so the new style
# reuses the table's existing border) and then paint it red/white.
# ---------------------------------------------------------------------

$borderTemplate = $ws.Range("H3")

$reportedCells = @("K12", "K38")
foreach ($ref in $reportedCells) {
    $borderTemplate.Copy()
    $ws.Range($ref).PasteSpecial(-4122)   # xlPasteFormats

    $cell = $ws.Range($ref)
    $cell.Value2 = "Reporté"
    $cell.Interior.Color = 204            # RGB(204,0,0) -> FFCC0000
    $cell.Font.Color = 16777215           # RGB(255,255,255) -> FFFFFFFF
    $cell.HorizontalAlignment = -4108     # xlCenter
    $cell.VerticalAlignment = -4108       # xlCenter
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Column widths: widen K (now holds the "Reporté" label) and shrink the
# hidden helper column L down to almost nothing.
# ---------------------------------------------------------------------

$ws.Columns.Item(11).ColumnWidth = 22.17
$ws.Columns.Item(12).ColumnWidth = -0.65
$ws.Columns.Item(12).Hidden = $true

# ---------------------------------------------------------------------
# Move the selection to the second newly-flagged "Reporté" cell.
# ---------------------------------------------------------------------

$ws.Range("K38").Select()
